# Replace the HYPERLINK() formulas in column F (rows 3-9, the C40 job
# postings) with plain-text URL strings, matching the inline-string format
# already used for the Apply_Link column throughout the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$urls = @{
    3 = "https://c40.bamboohr.com/careers/697"
    4 = "https://c40.bamboohr.com/careers/698"
    5 = "https://c40.bamboohr.com/careers/699"
    6 = "https://c40.bamboohr.com/careers/700"
    7 = "https://c40.bamboohr.com/careers/701"
    8 = "https://c40.bamboohr.com/careers/702"
    9 = "https://c40.bamboohr.com/careers/703"
}

foreach ($row in $urls.Keys) {
    # Column F is the 6th column (Apply_Link).
    $ws.Cells.Item($row, 6).Value = $urls[$row]
}
